$wb = $excel.ActiveWorkbook

# --- Summary sheet: build the new descriptive tables ---
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("A1").Value = "Sheet"
$summary.Range("B1").Value = "Description"
$summary.Range("A1:B1").Font.Bold = $true

$summary.Range("A2").Value = "financial"
$summary.Range("B2").Value = "Includes relevant metrics for the cost effectiveness of the system such as NPV, annual bill savings, and system upfront costs"

$summary.Range("A3").Value = "lcc breakdown"
$summary.Range("B3").Value = "Breaks out costs into respective components, such as utility energy costs and system CAPEX and O&M costs"

$summary.Range("A4").Value = "home"
$summary.Range("B4").Value = "Metrics relevant to home performance such as annual load and home comfort"

$summary.Range("A5").Value = "external"
$summary.Range("B5").Value = "Metrics relevant to utility and other stakeholders. These include grid impacts and environmental impacts."

$summary.Range("A6").Value = "Technologies (PV, Storage, HVAC, HPWH, ERWH)"
$summary.Range("B6").Value = "Metrics on system upgrades and technology performance, such as system capacity. Mainly useful for determining case"

$summary.Range("A9").Value = "Metric Type"
$summary.Range("B9").Value = "Description"
$summary.Range("A9:B9").Font.Bold = $true

$summary.Range("A10").Value = "Metadata"
$summary.Range("B10").Value = "Used to identify case"

$summary.Range("A11").Value = "Absolute"
$summary.Range("B11").Value = "Absolute value of run, does not use baseline comparison. "

$summary.Range("A12").Value = "Comparison"
$summary.Range("B12").Value = "Comparison to baseline. Example of difference is annual home load is an absolute metric (no comparison) while annual load reduction from an upgrade is a comparison metric (is old home load minus new home load)"

$summary.Range("A15").Value = "Metric Importance"
$summary.Range("B15").Value = "Description"
$summary.Range("A15:B15").Font.Bold = $true

$summary.Range("A16").Value = "Metadata for Case"
$summary.Range("B16").Value = "Used to identify case (includes system sizing as well as case name and baseline type)"

$summary.Range("A17").Value = "Internal Use Only"
$summary.Range("B17").Value = "Not actual costs but instead internal REopt values for optimization. Important for internal understanding of process but not valuable for external facing values"

$summary.Range("A18").Value = "Low Priority"
$summary.Range("B18").Value = "Technical metrics useful for understanding results, but low priority for external facing values"

$summary.Range("A19").Value = "Medium Priority"
$summary.Range("B19").Value = "Metrics which may be useful to present in some cases. In many cases these metrics are secondary or supporting metrics to other results"

$summary.Range("A20").Value = "Medium-High Priority"
$summary.Range("B20").Value = "Useful metrics which can be included in external facing values. "

$summary.Range("A21").Value = "High Priority"
$summary.Range("B21").Value = "Metrics external audience would most likely want to see. Provides good overview of what various stakeholders care about. Metrics to highlight. "

$summary.Range("B21").Select()

# --- Metrics sheet: fix the metrics generation bug ---
$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("F53").Value = "Low Priority"

$metrics.Activate()
$excel.ActiveWindow.ScrollRow = 1
$metrics.Range("A35:F35").Select()

$summary.Activate()
$summary.Range("B21").Select()
